# Auto commit at 2025-09-30 10:13:45.25
# Update Metrics values, change active sheet/tab selections to match the
# authored diff.

$wb = $excel.ActiveWorkbook

# --- Update Metrics sheet numeric values (B2:B13) -------------------------
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 431438.39
$metrics.Range("B3").Value = 348723.86
$metrics.Range("B4").Value = 136095.97
$metrics.Range("B5").Value = 17148
$metrics.Range("B6").Value = 4350689.3100000005
$metrics.Range("B7").Value = 3676251.3400000003
$metrics.Range("B8").Value = 1265461.6499999999
$metrics.Range("B9").Value = 168308
$metrics.Range("B10").Value = 32816013.100000001
$metrics.Range("B11").Value = 30951472.899999999
$metrics.Range("B12").Value = 11547170.560000001
$metrics.Range("B13").Value = 1265935

# --- Update sheet selections (cosmetic view-state) -------------------------
# Metrics: selection moves to E41 (and loses tabSelected once "today" is
# activated below).
$metrics.Range("E41").Select()

# Chargingdata: drop the pinned topLeftCell, keep selection at G7.
$chargingdata = $wb.Worksheets.Item("Chargingdata")
$chargingdata.Range("G7").Select()

# zgmysj: selection moves from P90 to O79 within the frozen pane.
$zgmysj = $wb.Worksheets.Item("zgmysj")
$zgmysj.Range("O79").Select()

# today: becomes the active/selected tab, selection moves to E6, and the
# pinned topLeftCell is cleared.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E6").Select()

# Recalculate so the "today" sheet's Metrics!-linked formulas (and the
# TODAY()-1 cell) refresh to reflect the new Metrics values.
$excel.CalculateFull()
